# Update the crypto "symbol list" snapshot (Price / Volume(1h) columns) to the
# newer scrape values. The Price column (D) holds numeric-looking values that
# the sheet stores as *text* (so things like trailing zeros, e.g. "0.8100" or
# "0.0001500", are preserved verbatim). Assigning a plain numeric string via
# .Value would make Excel auto-convert the cell to a Number (losing the exact
# text formatting), so each numeric-looking value is entered with a leading
# apostrophe to force Text entry, then the cell style is reset back to
# "Normal" so we don't leave a stray quote-prefixed number format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Leading apostrophe forces Excel to treat a numeric-looking string as
    # literal text instead of converting it to a Number.
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2"  "247.26"
Set-TextValue "D3"  "22.73"
Set-TextValue "D4"  "5.291"
Set-TextValue "D5"  "0.05725"
Set-TextValue "D7"  "0.8093"
Set-TextValue "D8"  "0.8697"
Set-TextValue "D9"  "0.1430"
Set-TextValue "D10" "0.07422"
Set-TextValue "D12" "0.03129"
Set-TextValue "D13" "0.09399"
Set-TextValue "D14" "4.001"
Set-TextValue "D15" "0.001578"
Set-TextValue "D16" "0.04823"

$ws.Range("E17").Value = "16OneONEWorstin24h"

Set-TextValue "D18" "0.006147"
Set-TextValue "D19" "0.005122"
Set-TextValue "D20" "0.0009984"
Set-TextValue "D21" "0.0001500"
Set-TextValue "D22" "3.736"
Set-TextValue "D23" "6.318"
Set-TextValue "D24" "2.185"
Set-TextValue "D25" "0.3280"

Set-TextValue "D41" "0.006765"
Set-TextValue "D42" "0.1066"
Set-TextValue "D43" "0.003201"
Set-TextValue "D44" "0.007493"
Set-TextValue "D45" "0.00005612"

Set-TextValue "D48" "0.1764"
$ws.Range("E48").Value = "47BOLOBOLO"

Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.01010"
